# Generate Report for Archive
#
# The localization status report is regenerated: the "Ready for handoff"
# status label becomes "In Translation" everywhere it is used (the
# Overview sheet's per-language status columns, plus each language
# sheet's own Status column), and the now-narrower status columns are
# resized to fit the shorter label.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status label ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Narrow the status columns to fit the shorter label ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
